$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data down
$ws.Rows.Item(2).Insert()

# Match the date-formatted style used in column A, clear style on B:E
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:E2").ClearFormats()

# Update all row values (row 2 is new, rows 3-19 are updated forecasts)
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 0.4235526809466261
$ws.Range("D2").Value = 2008
$ws.Range("E2").ClearContents()

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = -0.571807692896309
$ws.Range("D3").Value = 2009
$ws.Range("E3").ClearContents()

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = 0.3486139762224783
$ws.Range("D4").Value = 2010
$ws.Range("E4").ClearContents()

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = -0.1384957661262898
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 0.592211799485276

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 1.566479473280147
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 0.5295895589954247

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 0.7307568962936939
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 0.5955791956549161

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.818818812164257
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 0.9193568360546411

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 0.9180054319587239
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 0.9697679806505821

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 1.984684278296656
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 0.970573649360662

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 1.755995812646982
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 1.083941060573212

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 1.946965557828384
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 1.343559319682996

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 1.06432145354225
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 0.8791375467670726

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 1.361817904277696
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 1.226035857429442

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -4.352425014431304
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = -1.242805832929039

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = -1.761645650979182
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 3.46909983288044

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 5.20787683103745
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = 1.517423464826884

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -0.9008525709169546
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.2090788898015949

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = 0.2738544794132824
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 0.5362040463673612
